$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.560.90'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.121.13'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.22'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.74'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.51'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.479'
$ws.Range("E11").Value = '  -1.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000248'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.76'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.630.85'
$ws.Range("E15").Value = '  -0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.370.92'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.11'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.110.98'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.64'
$ws.Range("E19").Value = '  +2.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '492.01'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.701'
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.83'
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '83.94'
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.12'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.58'
$ws.Range("E26").Value = '  +4.79%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.90'
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("E29").Value = '  -2.43%  '
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.35'
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.114'
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0948'
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.81'
$ws.Range("E35").Value = '  -2.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.972'
$ws.Range("E36").Value = '  -2.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '46.94'
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("E38").Value = '  -3.72%  '
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.308'
$ws.Range("E40").Value = '  -1.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.47'
$ws.Range("E41").Value = '  -2.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '386.16'
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.804.84'
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.58'
$ws.Range("E44").Value = '  -8.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0352'
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.29'
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.00'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.72'
$ws.Range("E51").Value = '  -1.48%  '
